# "troca equipe por aluno" -- replace team/"time(s)" based copy with
# student/"aluno" based copy in a handful of paragraphs, and move the
# _GoBack bookmark from the old edit location to the new one.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: split a range into two runs by dropping a transient bookmark
# at the split point before an edit, then removing it afterwards. Word
# (and this runtime) only merges *adjacent* same-formatted runs when a
# Range.Text assignment touches both sides, so a temporary bookmark at
# the boundary keeps the two sides from re-merging once it is deleted.
# ---------------------------------------------------------------------

function Mark-Split($pos, $name) {
    $d.Bookmarks.Add($name, $d.Range($pos, $pos)) | Out-Null
}

# =======================================================================
# EDIT 1: "O professor irá cadastrar os times. Times são compostos de
# jogadores." -> "O professor irá cadastrar os jogadores." and the
# _GoBack bookmark now wraps this whole (new) sentence.
# =======================================================================

$r = $d.Content
$r.Find.Execute("O professor irá cadastrar os times. Times são compostos de jogadores.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s1 = $r.Start
$e1 = $r.End

# Mark the split between the two surviving runs ("...os " | "jogadores.")
$splitPos = $s1 + "O professor irá cadastrar os ".Length
Mark-Split $splitPos "SPLIT1"

# Delete "times. Times são compostos de " out of the middle.
$delStart = $s1 + "O professor irá cadastrar os ".Length
$delEnd = $s1 + "O professor irá cadastrar os times. Times são compostos de ".Length
$d.Range($delStart, $delEnd).Text = ""

$newEnd1 = $e1 - ($delEnd - $delStart)

$d.Bookmarks.Item("SPLIT1").Delete()

# Wrap the whole resulting sentence with the (re-homed) _GoBack bookmark.
$d.Bookmarks.Add("_GoBack", $d.Range($s1, $newEnd1)) | Out-Null

# =======================================================================
# EDIT 2: "Professor -> Sistema (Cadastro de Times)" ->
# "Professor -> Sistema (Cadastro de Alunos)" split as
# "Profess" | "or -> Sistema (Cadastro de Alunos" | ")"
# =======================================================================

$r = $d.Content
$r.Find.Execute("Professor -> Sistema (Cadastro de Times)", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s2 = $r.Start
$e2 = $r.End

Mark-Split ($s2 + "Profess".Length) "SPLIT2A"
Mark-Split ($e2 - 1) "SPLIT2B"

$full2 = $d.Range($s2, $e2).Text
$timesStart = $s2 + $full2.IndexOf("Times")
$timesEnd = $timesStart + "Times".Length
$d.Range($timesStart, $timesEnd).Text = "Alunos"

$d.Bookmarks.Item("SPLIT2A").Delete()
$d.Bookmarks.Item("SPLIT2B").Delete()

# =======================================================================
# EDIT 3: "Professor clica no ícone de Cadastrar Times" ->
# "Professor clica no ícone de Cadastrar Alunos" split as
# "Professor c" | "lica no ícone de Cadastrar Alunos"
# =======================================================================

$r = $d.Content
$r.Find.Execute("Professor clica no ícone de Cadastrar Times", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s3 = $r.Start
$e3 = $r.End

Mark-Split ($s3 + "Professor c".Length) "SPLIT3"

$full3 = $d.Range($s3, $e3).Text
$timesStart3 = $s3 + $full3.IndexOf("Times")
$timesEnd3 = $timesStart3 + "Times".Length
$d.Range($timesStart3, $timesEnd3).Text = "Alunos"

$d.Bookmarks.Item("SPLIT3").Delete()

# =======================================================================
# EDIT 4: "de cadastro de times" -> "de cadastro de alunos". This text
# is already its own run (preceded by "Sistema abre a tela " in a
# separate run); keep that boundary intact while swapping the word.
# =======================================================================

$r = $d.Content
$r.Find.Execute("Sistema abre a tela de cadastro de times", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s4 = $r.Start
$e4 = $r.End

Mark-Split ($s4 + "Sistema abre a tela ".Length) "SPLIT4"

$full4 = $d.Range($s4, $e4).Text
$timesStart4 = $s4 + $full4.IndexOf("times")
$timesEnd4 = $timesStart4 + "times".Length
$d.Range($timesStart4, $timesEnd4).Text = "alunos"

$d.Bookmarks.Item("SPLIT4").Delete()

# =======================================================================
# EDIT 5: "Sistema valida time." -> "Sistema valida os alunos." split as
# "Sistema valida os alunos" | "."
# =======================================================================

$r = $d.Content
$r.Find.Execute("Sistema valida time.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s5 = $r.Start
$e5 = $r.End

Mark-Split ($e5 - 1) "SPLIT5"

$d.Range($s5, $e5 - 1).Text = "Sistema valida os alunos"

$d.Bookmarks.Item("SPLIT5").Delete()

# =======================================================================
# EDIT 6: "Professor confirma o cadastro do time" ->
# "Professor confirma o cadastro " (single run, trailing space kept).
# =======================================================================

$d.Content.Find.Execute("Professor confirma o cadastro do time", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Professor confirma o cadastro ", 2) | Out-Null

# =======================================================================
# EDIT 7: "Sistema cadastra o time no banco." ->
# "Sistema cadastra no banco." split as "Sistema cadastra" | " no banco."
# =======================================================================

$r = $d.Content
$r.Find.Execute("Sistema cadastra o time no banco.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s7 = $r.Start
$e7 = $r.End

Mark-Split ($s7 + "Sistema cadastra".Length) "SPLIT7"

$full7 = $d.Range($s7, $e7).Text
$oldMid = " o time"
$midStart = $s7 + $full7.IndexOf($oldMid)
$midEnd = $midStart + $oldMid.Length
$d.Range($midStart, $midEnd).Text = ""

$d.Bookmarks.Item("SPLIT7").Delete()

# =======================================================================
# EDIT 8: "Aluno " + (old _GoBack bookmark, now removed above) + "joga"
# -> single merged run "Aluno joga". The bookmark already moved away in
# EDIT 1, so the two runs here just need to be coalesced into one; a
# straight Range.Text = "Aluno joga" is a no-op (content unchanged) and
# this runtime only re-normalises runs touched by an actual content
# change, so round-trip through a placeholder string first.
# =======================================================================

$r = $d.Content
$r.Find.Execute("Aluno joga", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s8 = $r.Start
$e8 = $r.End
$placeholder = "placeholder8"
$d.Range($s8, $e8).Text = $placeholder
$d.Range($s8, $s8 + $placeholder.Length).Text = "Aluno joga"
